# Apply BoM content update for rcbus-opl3 (KiBot regenerated BoM: 2025-10-09)
#
# The upstream KiBot run renumbered several reference designators and
# refreshed a handful of capacitor/resistor values, as well as bumping
# the "Date:" field on the report header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# --- Date: field (D5) -------------------------------------------------
# "2025-10-09" looks like a date, and a plain Value assignment would make
# Excel re-interpret it (and reformat the cell) as a date serial number.
# To keep it as plain text (matching the original shared-string cell) and
# to avoid touching the cell's existing style, compute the literal text in
# a scratch cell via a formula, then copy/paste-special (values only) it
# into the target cell.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""2025-10-09"""
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)  # xlPasteValues
$scratch.ClearContents()

# --- Capacitor rows (9-13) --------------------------------------------
$ws.Range("D9").Value  = "C5"
$ws.Range("E9").Value  = "68pF"
$ws.Range("D10").Value = "C3 C4"
$ws.Range("E10").Value = "2.7nF"
$ws.Range("D11").Value = "C2 C8 C11 C12 C13"
$ws.Range("D12").Value = "C1 C6 C7 C10"
$ws.Range("E12").Value = "10uF"
$ws.Range("D13").Value = "C9"

# --- Resistor rows (17-19) ---------------------------------------------
$ws.Range("D17").Value = "R1"
$ws.Range("D18").Value = "R2 R3"
$ws.Range("B19").Value = "Resistor Resistor, small symbol"
$ws.Range("C19").Value = "R"
$ws.Range("D19").Value = "R4 R5 R6 R7 R8 R9 R10 R11"

# --- IC reference designator swap (rows 20-23) --------------------------
$ws.Range("D20").Value = "U3"
$ws.Range("D21").Value = "U2"
$ws.Range("D22").Value = "U1"
$ws.Range("D23").Value = "U4"
